# Work Experience working code
# Adds a new "WorkExperienceAddMore_Success" worksheet (a copy of the existing
# "WorkExperience_Success" layout) right before "WorkExperience_Success", and
# refreshes sample data on a couple of sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Duplicate "WorkExperience_Success" to create "WorkExperienceAddMore_Success"
#    positioned immediately before it (i.e. right after "Skills_Success").
# ---------------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("WorkExperience_Success")
$srcIndex = $srcSheet.Index
$srcSheet.Copy($srcSheet)
$newSheet = $wb.Worksheets.Item($srcIndex)
$newSheet.Name = "WorkExperienceAddMore_Success"

# Populate the new sheet's sample row with its own data set.
$newSheet.Range("A2").Value = "Automotive"
$newSheet.Range("B2").Value = "Software Engineer"
$newSheet.Range("C2").Value = "Cerner Pvt Ltd"
$newSheet.Range("D2").Value = "Jan"
$newSheet.Range("E2").Value = "2013"
$newSheet.Range("F2").Value = "Feb"
$newSheet.Range("G2").Value = "2015"
$newSheet.Range("H2").Value = ""
$newSheet.Range("I2").Value = "update_success"
$newSheet.Range("G2").Select()

# ---------------------------------------------------------------------------
# 2. Update the original "WorkExperience_Success" sample data.
# ---------------------------------------------------------------------------
$weSheet = $wb.Worksheets.Item("WorkExperience_Success")
$weSheet.Range("A2").Value = "Construction Industry"
$weSheet.Range("B2").Value = "Civil Engineer"
$weSheet.Range("C2").Value = "Sakha Global"
$weSheet.Range("D2").Value = "Jan"
$weSheet.Range("E2").Value = "2012"
$weSheet.Range("F2").Value = "Feb"
$weSheet.Range("G2").Value = "2013"
$weSheet.Range("H2").Value = ""
$weSheet.Range("I2").Value = "update_success"

# ---------------------------------------------------------------------------
# 3. Update "Login_Success" sample credentials.
# ---------------------------------------------------------------------------
$loginSheet = $wb.Worksheets.Item("Login_Success")
$loginSheet.Range("A2").Value = "1234567890"
$loginSheet.Range("B2").Value = "Welcome@1234"
$loginSheet.Range("B2").Select()

# ---------------------------------------------------------------------------
# 4. Make "WorkExperience_Success" the active sheet / tab, with column B
#    scrolled to the left edge, matching the author's saved view state.
# ---------------------------------------------------------------------------
$weSheet.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$weSheet.Range("G2").Select()
